# Fills in the four "half-year summary" table cells (handicraft, math,
# music, torah) of the certificate with the teacher's comments.

$d = $word.ActiveDocument

$melacha = "במחצית זאת עשינו מלאכות לסוכה, התאמנו על סריגה ועשינו מלאכות יד ותכשיטים לתפארה.`nהייתה אוירה טובה בכיתה וכולן נהנו.`nשמעון אתה ילד נהדר, בהצלחה!"

$cheshbon = "במחצית זאת למדנו את תורת המיספרים, הכרנו כל מספר לעומק, והיתחלנו עם פעולות חשבון בסיסיות,חיבור וחיסור, התקדמנו הרבה עם הספר ""חושבים 1"".והתכוננו לקראת השנה החדשה בההכרה מלמעלה כל כפל וחילוק,`nשמעון אתה ילד נפלא,עלה והצלח!!"

$musica = "במחצית זאת למדנו על עולם המוזיקה, על התווים ועל רמות הקול, התעסקנו עם שירים על מעגל השנה, הייתה אוירה כיפית ונחמדה.`nשמעון אתה תלמיד מדהים!"

$torah = "במחצית זאת למדנו חומש בראשית, למדנו והתפעלנו מבריאת העולם, עקידת יצחק וכו....`nשמעון אתה תלמיד מצוין, בהצלחה!"

$texts = @($melacha, $cheshbon, $musica, $torah)

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $table = $d.Tables.Item($i)
    $cell = $table.Cell(1, 2)
    $r = $cell.Range
    # trim the trailing cell-mark/paragraph-mark character before setting text
    $r.End = $r.End - 1
    $r.Text = $texts[$i - 1]
}
